$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.353.68"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "3.119.75"
$ws.Range("E3").Value = "  -4.55%  "
$ws.Range("E4").Value = "  -0.45%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.39%  "
$ws.Range("E7").Value = "  -0.62%  "
$ws.Range("D8").Value = "3.117.74"
$ws.Range("E8").Value = "  -1.03%  "
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.96"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.28%  "
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000245"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.34%  "
$ws.Range("D15").Value = "3.628.10"
$ws.Range("E15").Value = "  -4.63%  "
$ws.Range("E16").Value = "  -1.65%  "
$ws.Range("E17").Value = "  +2.38%  "
$ws.Range("D18").Value = "64.064.43"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "3.110.77"
$ws.Range("E19").Value = "  -2.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "470.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.740"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.52%  "
$ws.Range("B24").Value = "Fetch.AI"
$ws.Range("C24").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.49%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.11%  "
$ws.Range("E30").Value = "  +1.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("E33").Value = "  +6.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.92%  "
$ws.Range("D35").Value = "0.0₃0856"
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("E36").Value = "  +1.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.18"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.56%  "
$ws.Range("E39").Value = "  -2.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "458.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "51.03"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.79%  "
$ws.Range("E43").Value = "  +2.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0372"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "2.856.99"
$ws.Range("E45").Value = "  -2.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.110"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.90%  "
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.49%  "
$ws.Range("B51").Value = "USDe"
$ws.Range("C51").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.04%  "
